$d = $word.ActiveDocument

# --- Define the three new character styles (appended at end of styles.xml) ---

$styleGaNStyle = $d.Styles.Add("GaNStyle", 2)
$styleGaNStyle.Font.Name = "Calibri"
$styleGaNStyle.Font.Size = 14

$styleGaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$styleGaNParagraph.Font.Name = "Calibri"
$styleGaNParagraph.Font.Size = 10

$styleGaNLinks = $d.Styles.Add("GaNLinks", 2)
$styleGaNLinks.Font.Name = "Calibri"
$styleGaNLinks.Font.Size = 9.5
$styleGaNLinks.Font.Bold = $true
$styleGaNLinks.Font.Color = 8388608
$styleGaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every occurrence of the campaign-dates paragraph ---

$searchDates = "2022: Daty kampanii używające Gwiazdozbiór Oriona: 16-25 stycznia, 14-23 lutego, 14-24 marca"
$range = $d.Content
$range.Start = 0
$range.Find.ClearFormatting()
$found = $range.Find.Execute($searchDates, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
    $found = $range.Find.Execute($searchDates, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Uczestniczysz w ogólnoświatowym..." paragraph ---

$searchParagraph = "Uczestniczysz w ogólnoświatowym przedsięwzięciu, którego celem jest obserwacja i odnotowanie najsłabszych widocznych gwiazd w celu zmierzenia zanieczyszczenia światłem w danym miejscu. Poprzez zlokalizowanie i obserwację  Gwiazdozbiór Oriona na nocnym niebie oraz porównanie go do map nieba ludzie z całego świata będą mogli dowiedzieć się jaki wkład światło emitowane przez ich społeczność wnosi do  zanieczyszczenia światłem. To co dodasz do internetowej bazy danych pomoże udokumentować widoczne nocne niebo."
$range2 = $d.Content
$range2.Start = 0
$range2.Find.ClearFormatting()
$found2 = $range2.Find.Execute($searchParagraph, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $range2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Jenika Hollana, CzechGlobe..." run ---

$searchLinks = " Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$range3 = $d.Content
$range3.Start = 0
$range3.Find.ClearFormatting()
$found3 = $range3.Find.Execute($searchLinks, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $range3.Style = "GaNLinks"
}
